$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 385; this pushes the existing
# rows 385:476 down to 386:477 and extends the sheet dimension to A1:R477.
$ws.Rows("385:385").Insert()

# Populate the newly inserted row 385 with the new weekly price record.
$ws.Range("A385").Value = 10
$ws.Range("B385").Value = "Vega Modelo de Temuco"
$ws.Range("C385").Value = "La Araucanía"
$ws.Range("D385").Value = 44855
$ws.Range("E385").Value = 9
$ws.Range("F385").Value = 100114014
$ws.Range("G385").Value = "Betarraga"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 100
$ws.Range("K385").Value = 12000
$ws.Range("L385").Value = 12000
$ws.Range("M385").Value = 12000
$ws.Range("N385").Value = "$/docena de paquetes"
$ws.Range("O385").Value = "Región del Maule"
$ws.Range("P385").Value = 1000
$ws.Range("Q385").Value = 12
$ws.Range("R385").Value = "Hortaliza"
